$d = $word.ActiveDocument

# Update the date/title line
$d.Content.Find.Execute("2025-12-07 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-08 Monday", 2) | Out-Null

# Update the division problems in the table, addressed by row/column so that
# cells whose new value equals another cell's old value are not double-replaced.
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "87÷4=21, 3"  # was "53÷7=7, 4"
$t.Cell(1, 2).Range.Text = "55÷5=11, 0"  # was "54÷2=27, 0"
$t.Cell(1, 3).Range.Text = "83÷4=20, 3"  # was "98÷3=32, 2"
$t.Cell(1, 4).Range.Text = "85÷3=28, 1"  # was "69÷7=9, 6"
$t.Cell(1, 5).Range.Text = "70÷8=8, 6"  # was "89÷2=44, 1"
$t.Cell(5, 1).Range.Text = "89÷8=11, 1"  # was "90÷7=12, 6"
$t.Cell(5, 2).Range.Text = "70÷6=11, 4"  # was "89÷6=14, 5"
$t.Cell(5, 3).Range.Text = "74÷6=12, 2"  # was "87÷5=17, 2"
$t.Cell(5, 4).Range.Text = "61÷8=7, 5"  # was "90÷2=45, 0"
$t.Cell(5, 5).Range.Text = "41÷8=5, 1"  # was "39÷5=7, 4"
$t.Cell(9, 1).Range.Text = "17÷7=2, 3"  # was "61÷5=12, 1"
$t.Cell(9, 2).Range.Text = "31÷6=5, 1"  # was "52÷4=13, 0"
$t.Cell(9, 3).Range.Text = "58÷6=9, 4"  # was "66÷8=8, 2"
$t.Cell(9, 4).Range.Text = "40÷2=20, 0"  # was "52÷2=26, 0"
$t.Cell(9, 5).Range.Text = "57÷6=9, 3"  # was "58÷7=8, 2"
$t.Cell(13, 1).Range.Text = "68÷9=7, 5"  # was "79÷7=11, 2"
$t.Cell(13, 2).Range.Text = "54÷6=9, 0"  # was "77÷2=38, 1"
$t.Cell(13, 3).Range.Text = "64÷2=32, 0"  # was "68÷8=8, 4"
$t.Cell(13, 4).Range.Text = "98÷2=49, 0"  # was "74÷9=8, 2"
$t.Cell(13, 5).Range.Text = "75÷4=18, 3"  # was "70÷9=7, 7"
$t.Cell(17, 1).Range.Text = "15÷7=2, 1"  # was "83÷4=20, 3"
$t.Cell(17, 2).Range.Text = "46÷2=23, 0"  # was "19÷3=6, 1"
$t.Cell(17, 3).Range.Text = "54÷7=7, 5"  # was "92÷6=15, 2"
$t.Cell(17, 4).Range.Text = "12÷8=1, 4"  # was "35÷6=5, 5"
$t.Cell(17, 5).Range.Text = "94÷6=15, 4"  # was "66÷4=16, 2"
